$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Content changes: new IG canonical base + regeneration date -------------

# Metadata sheet: StructureDefinition URL (row 2) and generation Date (row 8)
$meta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/social-activity"
$meta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# Elements sheet: the Extension.url fixed value (R5) and the
# Binding Value Set URL (Z6) both carry the old canonical base URL - update both
# so every occurrence of the shared string is consistent with the new IG home.
$elements.Range("R5").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/social-activity"
$elements.Range("Z6").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/ValueSet/social-activity-vs"

# --- Column widths: the regenerated output re-ran "best fit" sizing on the --
# --- Elements sheet, shrinking every auto-fit column --------------------------

$elements.Columns("A:A").ColumnWidth = 15.666666666666666
$elements.Columns("B:B").ColumnWidth = 15.666666666666666
$elements.Columns("C:C").ColumnWidth = 9.0
$elements.Columns("D:D").ColumnWidth = 6.166666666666667
$elements.Columns("E:E").ColumnWidth = 4.5
$elements.Columns("F:F").ColumnWidth = 3.1666666666666665
$elements.Columns("G:G").ColumnWidth = 3.5
$elements.Columns("H:H").ColumnWidth = 11.833333333333334
$elements.Columns("I:I").ColumnWidth = 9.666666666666666
$elements.Columns("J:J").ColumnWidth = 19.833333333333332
$elements.Columns("K:K").ColumnWidth = 13.5
$elements.Columns("L:L").ColumnWidth = 99.83333333333333
$elements.Columns("M:M").ColumnWidth = 99.83333333333333
$elements.Columns("N:N").ColumnWidth = 99.83333333333333
$elements.Columns("O:O").ColumnWidth = 11.5
$elements.Columns("P:P").ColumnWidth = 19.833333333333332
$elements.Columns("Q:Q").ColumnWidth = 19.833333333333332
$elements.Columns("R:R").ColumnWidth = 19.833333333333332
$elements.Columns("S:S").ColumnWidth = 19.833333333333332
$elements.Columns("T:T").ColumnWidth = 7.0
$elements.Columns("U:U").ColumnWidth = 12.833333333333334
$elements.Columns("V:V").ColumnWidth = 13.166666666666666
$elements.Columns("W:W").ColumnWidth = 14.166666666666666
$elements.Columns("X:X").ColumnWidth = 13.833333333333334
$elements.Columns("Y:Y").ColumnWidth = 16.166666666666668
$elements.Columns("Z:Z").ColumnWidth = 54.166666666666664
$elements.Columns("AA:AA").ColumnWidth = 4.166666666666667
$elements.Columns("AB:AB").ColumnWidth = 17.166666666666668
$elements.Columns("AC:AC").ColumnWidth = 33.666666666666664
$elements.Columns("AD:AD").ColumnWidth = 12.666666666666666
$elements.Columns("AE:AE").ColumnWidth = 10.5
$elements.Columns("AF:AF").ColumnWidth = 14.166666666666666
$elements.Columns("AG:AG").ColumnWidth = 7.333333333333333
$elements.Columns("AH:AH").ColumnWidth = 7.666666666666667
$elements.Columns("AI:AI").ColumnWidth = 99.83333333333333
$elements.Columns("AK:AK").ColumnWidth = 18.666666666666668

# Re-assert the hidden columns (changing ColumnWidth can reset visibility)
$elements.Columns("C:C").Hidden = $true
$elements.Columns("D:D").Hidden = $true
$elements.Columns("AE:AE").Hidden = $true
$elements.Columns("AF:AF").Hidden = $true
$elements.Columns("AG:AG").Hidden = $true
